$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1223.1428
$ws.Range("I43").Value = 1123.5
$ws.Range("J43").Value = 1297.875
$ws.Range("K43").Value = 1123.5
$ws.Range("L43").Value = 1297.875
$ws.Range("M43").Value = -1054.5
$ws.Range("N43").Value = -1435.875
# Row 62
$ws.Range("H62").Value = 3249.8845
$ws.Range("I62").Value = 2878.4375
$ws.Range("J62").Value = 3844.2
$ws.Range("K62").Value = 2878.4375
$ws.Range("L62").Value = 3844.2
$ws.Range("M62").Value = -2254.4375
$ws.Range("N62").Value = -5092.2
# Row 65
$ws.Range("H65").Value = 3249.8845
$ws.Range("I65").Value = 2878.4375
$ws.Range("J65").Value = 3844.2
$ws.Range("K65").Value = 14392.1875
$ws.Range("L65").Value = 19221
$ws.Range("M65").Value = -11272.1875
$ws.Range("N65").Value = -25461
# Row 130
$ws.Range("H130").Value = 11997
$ws.Range("J130").Value = 11997
$ws.Range("L130").Value = 11997
$ws.Range("N130").Value = -22037
# Row 132
$ws.Range("H132").Value = 16695.467
$ws.Range("I132").Value = 16695.467
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 50086.401
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -47556.401
$ws.Range("N132").ClearContents()
# Row 135
$ws.Range("H135").Value = 569.2683
$ws.Range("I135").Value = 569.2683
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5123.414699999999
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2588.414699999999
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 1439.36
$ws.Range("I137").Value = 1129.4783
$ws.Range("J137").Value = 5003
$ws.Range("K137").Value = 3388.4349
$ws.Range("L137").Value = 15009
$ws.Range("M137").Value = -838.4349000000002
$ws.Range("N137").Value = -20109
# Row 138
$ws.Range("H138").Value = 2596.8904
$ws.Range("I138").Value = 1258.3704
$ws.Range("J138").Value = 6401.1055
$ws.Range("K138").Value = 3775.1112
$ws.Range("L138").Value = 19203.3165
$ws.Range("M138").Value = 1364.8888
$ws.Range("N138").Value = -29483.3165
# Row 141
$ws.Range("H141").Value = 7574.345
$ws.Range("I141").Value = 1134.3334
$ws.Range("J141").Value = 19776.475
$ws.Range("K141").Value = 3403.0002
$ws.Range("L141").Value = 59329.425
$ws.Range("M141").Value = 1776.9998
$ws.Range("N141").Value = -69689.42499999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2313.71
$ws.Range("I32").Value = 2078.1133
$ws.Range("J32").Value = 9931.333
$ws.Range("K32").Value = 2078.1133
$ws.Range("L32").Value = 9931.333
$ws.Range("M32").Value = -1791.1133
$ws.Range("N32").Value = -10505.333
# Row 63
$ws.Range("H63").Value = 5366.5557
$ws.Range("I63").Value = 6583.1665
$ws.Range("K63").Value = 6583.1665
$ws.Range("M63").Value = -5897.1665
# Row 66
$ws.Range("H66").Value = 5366.5557
$ws.Range("I66").Value = 6583.1665
$ws.Range("K66").Value = 32915.8325
$ws.Range("M66").Value = -29483.8325
# Row 101
$ws.Range("H101").Value = 47999
$ws.Range("J101").Value = 47999
$ws.Range("L101").Value = 47999
$ws.Range("N101").Value = -54489
# Row 124
$ws.Range("H124").Value = 35472
$ws.Range("J124").Value = 35472
$ws.Range("L124").Value = 35472
$ws.Range("N124").Value = -45292

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 919.7692
$ws.Range("I107").Value = 897.4
$ws.Range("J107").Value = 994.3333
$ws.Range("K107").Value = 897.4
$ws.Range("L107").Value = 994.3333
$ws.Range("M107").Value = 1022.6
$ws.Range("N107").Value = -4834.3333
# Row 134
$ws.Range("H134").Value = 1653.2632
$ws.Range("I134").Value = 1199.8667
$ws.Range("J134").Value = 3353.5
$ws.Range("K134").Value = 3599.6001
$ws.Range("L134").Value = 10060.5
$ws.Range("M134").Value = -1064.6001
$ws.Range("N134").Value = -15130.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2434.7446
$ws.Range("I31").Value = 1942.1818
$ws.Range("J31").Value = 2868.2
$ws.Range("K31").Value = 1942.1818
$ws.Range("L31").Value = 2868.2
$ws.Range("M31").Value = -1647.1818
$ws.Range("N31").Value = -3458.2
# Row 34
$ws.Range("H34").Value = 2434.7446
$ws.Range("I34").Value = 1942.1818
$ws.Range("J34").Value = 2868.2
$ws.Range("K34").Value = 1942.1818
$ws.Range("L34").Value = 2868.2
$ws.Range("M34").Value = -1740.1818
$ws.Range("N34").Value = -3272.2
# Row 107
$ws.Range("H107").Value = 611.2143
$ws.Range("I107").Value = 285
$ws.Range("J107").Value = 893.93335
$ws.Range("K107").Value = 285
$ws.Range("L107").Value = 893.93335
$ws.Range("M107").Value = 1635
$ws.Range("N107").Value = -4733.93335
# Row 132
$ws.Range("H132").Value = 1564.3077
$ws.Range("I132").Value = 568.5641
$ws.Range("J132").Value = 4551.5386
$ws.Range("K132").Value = 1705.6923
$ws.Range("L132").Value = 13654.6158
$ws.Range("M132").Value = 824.3076999999998
$ws.Range("N132").Value = -18714.6158

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 5925
$ws.Range("J131").Value = 6755.263
$ws.Range("L131").Value = 20265.789
$ws.Range("N131").Value = -30345.789

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1623
$ws.Range("I2").Value = 2051.8
$ws.Range("J2").Value = 1355
$ws.Range("K2").Value = 2051.8
$ws.Range("L2").Value = 1355
$ws.Range("M2").Value = -1938.8
$ws.Range("N2").Value = -1581
# Row 69
$ws.Range("H69").Value = 26000
$ws.Range("J69").Value = 26000
$ws.Range("L69").Value = 26000
$ws.Range("N69").Value = -27498
# Row 72
$ws.Range("H72").Value = 26000
$ws.Range("J72").Value = 26000
$ws.Range("L72").Value = 78000
$ws.Range("N72").Value = -85488
# Row 107
$ws.Range("H107").Value = 1368.75
$ws.Range("I107").Value = 1796.6666
$ws.Range("J107").Value = 85
$ws.Range("K107").Value = 1796.6666
$ws.Range("L107").Value = 85
$ws.Range("M107").Value = 123.3334
$ws.Range("N107").Value = -3925
# Row 126
$ws.Range("H126").Value = 1900.4445
$ws.Range("I126").Value = 1701.1818
$ws.Range("J126").Value = 2213.5715
$ws.Range("K126").Value = 5103.5454
$ws.Range("L126").Value = 6640.7145
$ws.Range("M126").Value = -2633.5454
$ws.Range("N126").Value = -11580.7145

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 12347871
$ws.Range("I136").Value = 2898.5
$ws.Range("J136").Value = 37037816
$ws.Range("K136").Value = 8695.5
$ws.Range("L136").Value = 111113448
$ws.Range("M136").Value = -6145.5
$ws.Range("N136").Value = -111118548
# Row 139
$ws.Range("H139").Value = 35881.316
$ws.Range("J139").Value = 36763.61
$ws.Range("L139").Value = 36763.61
$ws.Range("N139").Value = -47043.61

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 374.8125
$ws.Range("I107").Value = 305.7143
$ws.Range("J107").Value = 428.55554
$ws.Range("K107").Value = 917.1428999999999
$ws.Range("L107").Value = 1285.66662
$ws.Range("M107").Value = 1002.8571
$ws.Range("N107").Value = -5125.66662
# Row 136
$ws.Range("H136").Value = 1493.7407
$ws.Range("I136").Value = 777.7143
$ws.Range("K136").Value = 2333.1429
$ws.Range("M136").Value = 216.8571000000002
